$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate()

# CreateDate (row 18) type changes from DATE to TIMESTAMP, length cleared
$ws.Range("D18").Value = "TIMESTAMP"
$ws.Range("E18").Value = $null

# LastUpdate (row 20) type changes from DATE to TIMESTAMP, length cleared
$ws.Range("D20").Value = "TIMESTAMP"
$ws.Range("E20").Value = $null

# Update selection to G19 on the DBD sheet
$ws.Range("G19").Select()
